# Fittings_level_tmp.xlsx - "woah more fittings stuff cool"
#
# Semantic changes applied:
#   - Hide the two leftmost helper columns (A:B)
#   - Hide the blank spacer row at the top of the sheet (row 1)
#   - Move the active selection/scroll position to reflect where the
#     user ended up working (top-left near C2, active cell E26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse and hide helper columns A:B
$ws.Columns("A:B").ColumnWidth = 0
$ws.Columns("A:B").Hidden = $true

# Hide the thin spacer row above the header
$ws.Rows("1:1").Hidden = $true

# Scroll so row/column C2 sits at the top-left of the viewport
$ws.Range("C2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 3
try {
  $excel.ActiveWindow.TopLeftCell = $ws.Range("C2")
} catch {
}

# Restore the active cell/selection to where the user left off
$ws.Range("E26").Select()
